$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# Header for new hole_id column - copy style (bold/border/alignment) from B1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "hole_id"

# hole_id values for rows 2-32, in order
$holeIds = @(
    "BRG_16_01",
    "BRG_01_06",
    "BRG_01_01",
    "BRG_16_02",
    "BRG_01_03",
    "BRG_05_15",
    "BRG_05_14",
    "BRG_16_09",
    "BRG_05_02",
    "ECO_09_05",
    "BRG_05_13",
    "BRG_16_04B",
    "BRG_01_02",
    "BRG_05_03",
    "BRG_05_05",
    "BRG_16_05",
    "BRG_16_07",
    "BRG_16_08",
    "BRG_08_01",
    "BRG_05_12",
    "BRG_05_10",
    "BRG_05_04",
    "BRG_01_08",
    "ECO_09_03",
    "BRG_16_06",
    "BRG_16_04A",
    "BRG_01_09",
    "BRG_05_01",
    "BRG_05_11",
    "BRG_05_09",
    "ECO_09_02"
)

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
